$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Dear Valued Customer, " (the
# greeting is split across two runs: the bold/yellow-highlighted "Dear Valued
# Customer," run and a trailing bold space run). We replace the whole
# paragraph's text with the new greeting as completely unformatted text.

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "Dear Valued Customer*") {
        $r = $para.Range
        # Exclude the trailing paragraph mark from the range we rewrite.
        $r.End = $r.End - 1

        # First clear the range's text completely: an empty range picks up
        # no run-level character formatting (rPr) of its own, only the
        # paragraph mark's formatting remains on the pPr/rPr.
        $r.Text = ""
        # Now insert the replacement text; the newly typed run comes out
        # with no rPr at all (no bold, no highlight), matching the target.
        $r.Text = "Dear APPLIED MEDICAL RESOURCES,"

        $found = $true
        break
    }
}

Write-Output "replaced: $found"
